{"js": "// Update the worksheet date and the 25 division problems/answers in the\n// table to the new day's content.\n//\n// The table's populated rows are every 4th row (rows 0, 4, 8, 12, 16 \u2014\n// interleaved with 3 blank filler rows), each holding 5 cells. We target\n// cells by (row, column) structurally instead of matching on the old\n// text, so values that happen to collide with another cell's *new* text\n// (e.g. \"42\u00f72=21, 0\" is the new value of two different cells, and\n// \"20\u00f78=2, 4\" is simultaneously an old value in one cell and a new value\n// in another) are never mismatched.\nconst newDate = \"2025-12-23 Tuesday\";\n\n// One row per table \"block\"; values left-to-right match the table's\n// column order.\nconst newRowValues = [\n  [\"46\u00f75=9, 1\", \"38\u00f73=12, 2\", \"93\u00f73=31, 0\", \"30\u00f72=15, 0\", \"32\u00f77=4, 4\"],\n  [\"55\u00f78=6, 7\", \"37\u00f79=4, 1\", \"71\u00f73=23, 2\", \"42\u00f72=21, 0\", \"14\u00f78=1, 6\"],\n  [\"42\u00f72=21, 0\", \"33\u00f76=5, 3\", \"63\u00f76=10, 3\", \"14\u00f77=2, 0\", \"82\u00f73=27, 1\"],\n  [\"84\u00f75=16, 4\", \"94\u00f73=31, 1\", \"23\u00f74=5, 3\", \"29\u00f75=5, 4\", \"20\u00f78=2, 4\"],\n  [\"16\u00f72=8, 0\", \"94\u00f76=15, 4\", \"30\u00f78=3, 6\", \"38\u00f79=4, 2\", \"22\u00f72=11, 0\"],\n];\n\n// Update the centered date line at the top of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nif (paragraphs.items.length === 0) {\n  throw new Error(\"Document has no paragraphs; cannot update date line.\");\n}\nparagraphs.items[0].insertText(newDate, Word.InsertLocation.replace);\n\n// Update the table of division problems.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table of division problems but found none.\");\n}\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst dataRowIndices = [0, 4, 8, 12, 16];\nif (rows.items.length < 17) {\n  throw new Error(`Expected at least 17 table rows, found ${rows.items.length}`);\n}\n\nfor (let i = 0; i < dataRowIndices.length; i++) {\n  const row = rows.items[dataRowIndices[i]];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const values = newRowValues[i];\n  if (cells.items.length < values.length) {\n    throw new Error(\n      `Row ${dataRowIndices[i]} has ${cells.items.length} cells, expected ${values.length}`\n    );\n  }\n\n  // Replace text on the cell's existing paragraph (not `cell.body`) so the\n  // paragraph/run formatting (left alignment, TimeNewRoman font, size 30)\n  // already on that run is preserved instead of being reset to defaults.\n  for (let c = 0; c < values.length; c++) {\n    const cellParagraphs = cells.items[c].body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].insertText(values[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division problems/answers in the\n# table to the new day's content.\n#\n# The table's populated rows are every 4th row (1-based rows 1, 5, 9, 13,\n# 17 - interleaved with 3 blank filler rows), each holding 5 cells. Cells\n# are targeted by (row, column) structurally instead of matching on the\n# old text, so values that happen to collide with another cell's *new*\n# text (e.g. \"42\u00f72=21, 0\" is the new value of two different cells, and\n# \"20\u00f78=2, 4\" is simultaneously an old value in one cell and a new value\n# in another) are never mismatched. Assigning directly to `Range.Text`\n# keeps the existing run/paragraph formatting (left alignment,\n# TimeNewRoman, size 30) on the cell, it only swaps the text.\n\n$d = $word.ActiveDocument\n\n# Update the centered date line at the top of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2025-12-23 Tuesday\"\n\n# Update the table of division problems.\n$table = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$newRowValues = @(\n  @(\"46\u00f75=9, 1\", \"38\u00f73=12, 2\", \"93\u00f73=31, 0\", \"30\u00f72=15, 0\", \"32\u00f77=4, 4\"),\n  @(\"55\u00f78=6, 7\", \"37\u00f79=4, 1\", \"71\u00f73=23, 2\", \"42\u00f72=21, 0\", \"14\u00f78=1, 6\"),\n  @(\"42\u00f72=21, 0\", \"33\u00f76=5, 3\", \"63\u00f76=10, 3\", \"14\u00f77=2, 0\", \"82\u00f73=27, 1\"),\n  @(\"84\u00f75=16, 4\", \"94\u00f73=31, 1\", \"23\u00f74=5, 3\", \"29\u00f75=5, 4\", \"20\u00f78=2, 4\"),\n  @(\"16\u00f72=8, 0\", \"94\u00f76=15, 4\", \"30\u00f78=3, 6\", \"38\u00f79=4, 2\", \"22\u00f72=11, 0\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $rowIndex = $dataRows[$i]\n  $values = $newRowValues[$i]\n  for ($c = 0; $c -lt $values.Length; $c++) {\n    $cell = $table.Cell($rowIndex, $c + 1)\n    $cell.Range.Text = $values[$c]\n  }\n}\n"}
